# Lotomania workbook update: append the new lottery draw results
# (contests 2808-2823) below the existing history table, matching the
# author's latest export of LoteriasExcel/Lotomania_edt.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new draw rows (contest number in col A, the 20 drawn numbers in B:U) ---
$newDraws = @(
    @(2808,3,11,12,15,25,26,27,31,32,39,52,57,58,65,72,74,80,83,84,97),
    @(2809,8,17,19,23,25,30,41,48,50,53,65,71,72,75,81,84,88,95,96,99),
    @(2810,3,4,6,8,20,22,28,33,37,44,53,54,55,58,61,63,73,77,82,95),
    @(2811,5,22,28,33,37,44,46,47,48,54,55,57,68,70,77,79,81,87,92,96),
    @(2812,2,6,7,17,20,22,25,32,37,43,49,55,58,71,74,79,81,83,87,95),
    @(2813,7,9,16,20,36,40,47,50,54,58,66,74,76,77,80,84,86,89,90,99),
    @(2814,2,5,18,19,22,25,26,28,31,35,40,43,65,71,72,73,79,81,93,96),
    @(2815,1,4,5,8,9,32,38,43,52,53,58,62,63,64,68,78,79,86,88,89),
    @(2816,8,13,14,21,26,28,31,34,43,44,50,52,54,64,73,74,78,84,90,91),
    @(2817,0,1,4,9,16,20,26,35,49,54,65,72,73,84,85,87,93,94,97,98),
    @(2818,3,4,8,15,30,31,32,34,52,53,55,57,58,62,71,72,79,86,89,93),
    @(2819,0,13,16,19,24,28,35,38,44,46,52,53,62,69,71,75,77,79,86,89),
    @(2820,0,1,2,16,18,22,29,30,35,37,41,47,50,65,67,70,75,86,93,97),
    @(2821,8,16,17,19,23,25,30,36,40,44,48,51,64,66,77,83,88,94,95,98),
    @(2822,12,13,18,20,29,34,37,38,41,50,52,58,60,70,72,74,76,89,90,99),
    @(2823,2,8,9,12,13,17,18,26,41,44,61,67,76,77,79,80,82,84,87,90)
)

$firstRow = 302
$r = $firstRow
foreach ($draw in $newDraws) {
    for ($i = 0; $i -lt $draw.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $draw[$i]
    }
    # column V has no data for these rows, just shares the new row style
    $ws.Cells.Item($r, 22).Value = ""
    $r++
}
$lastRow = $r - 1

# a trailing, otherwise-empty row whose only populated cell is V
$trailingRow = $lastRow + 1
$ws.Cells.Item($trailingRow, 22).Value = ""

# Stamp every new cell (A:V across the draw rows, plus the lone V cell on
# the trailing row) with the same plain/default formatting so they pick up
# a fresh style record rather than inheriting stray formatting.
$ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, 22)).WrapText = $false
$ws.Cells.Item($trailingRow, 22).WrapText = $false

# --- move the viewport / selection down to the new bottom of the table ---
$ws.Application.Goto($ws.Range("V" + ($trailingRow + 4)))
